$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.75
$ws.Range("I2").Value = 1.85
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 2.4
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.49
$ws.Range("R2").Value = 2.45
$ws.Range("S2").Value = 1.88
$ws.Range("T2").Value = 2.02
$ws.Range("U2").Value = 2.3
$ws.Range("V2").Value = 1.62
$ws.Range("W2").Value = 1.29
$ws.Range("X2").Value = 3.5
$ws.Range("Y2").Value = 1.53
$ws.Range("Z2").Value = 2.38
$ws.Range("AF2").Value = 29
$ws.Range("AI2").Value = 13
$ws.Range("AK2").Value = 126
$ws.Range("AL2").Value = 12
$ws.Range("AM2").Value = 12
$ws.Range("AO2").Value = 19
$ws.Range("AQ2").Value = 21

# Row 3
$ws.Range("G3").Value = 1.9
$ws.Range("I3").Value = 4.33
$ws.Range("J3").Value = 2.62
$ws.Range("K3").Value = 1.92
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.47
$ws.Range("P3").Value = 2.5
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 1.24
$ws.Range("V3").Value = 1.13
$ws.Range("AB3").Value = 7.5
$ws.Range("AD3").Value = 15
$ws.Range("AM3").Value = 21
$ws.Range("AO3").Value = 51

# Row 4
$ws.Range("G4").Value = 1.73
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.2
$ws.Range("V4").Value = 1.22
$ws.Range("AB4").Value = 7
$ws.Range("AC4").Value = 8.5
$ws.Range("AE4").Value = 15

# Row 5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 4.1
$ws.Range("I5").Value = 1.57
$ws.Range("K5").Value = 2.5
$ws.Range("L5").Value = 2.1
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 17
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 5
$ws.Range("Q5").Value = 1.57
$ws.Range("R5").Value = 2.35
$ws.Range("S5").Value = 1.87
$ws.Range("T5").Value = 1.87
$ws.Range("U5").Value = 2.37
$ws.Range("V5").Value = 1.5
$ws.Range("W5").Value = 1.29
$ws.Range("X5").Value = 3.5
$ws.Range("Y5").Value = 1.62
$ws.Range("Z5").Value = 2.2
$ws.Range("AA5").Value = 19
$ws.Range("AC5").Value = 17
$ws.Range("AF5").Value = 34
$ws.Range("AG5").Value = 17
$ws.Range("AH5").Value = 8.5
$ws.Range("AL5").Value = 9.5
$ws.Range("AM5").Value = 9

# Row 7
$ws.Range("G7").Value = 1.48
$ws.Range("H7").Value = 4.75
$ws.Range("I7").Value = 6.25
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 2.4
$ws.Range("L7").Value = 6.5
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 15
$ws.Range("Q7").Value = 1.67
$ws.Range("R7").Value = 2.1
$ws.Range("Y7").Value = 1.91
$ws.Range("Z7").Value = 1.91
$ws.Range("AB7").Value = 7
$ws.Range("AD7").Value = 10
$ws.Range("AF7").Value = 26
$ws.Range("AH7").Value = 8.5
$ws.Range("AI7").Value = 19
$ws.Range("AK7").Value = 301
$ws.Range("AL7").Value = 17
$ws.Range("AM7").Value = 34
$ws.Range("AN7").Value = 19
$ws.Range("AO7").Value = 67

